$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 49, shifting old rows 49->50 and 50->51 down.
$ws.Rows.Item(49).Insert()

# Populate the new row 49 with values (same A,B,C,E,F,G,H,I,N,Q,R as old row 49;
# new D,J,K,L,M,O,P values per the target diff).
$ws.Cells.Item(49, 1).Value = 5
$ws.Cells.Item(49, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(49, 3).Value = 'Maule'
$ws.Cells.Item(49, 4).Value = 44516
$ws.Cells.Item(49, 5).Value = 7
$ws.Cells.Item(49, 6).Value = 300000000
$ws.Cells.Item(49, 7).Value = 'Espárragos'
$ws.Cells.Item(49, 8).Value = 'Verde'
$ws.Cells.Item(49, 9).Value = 'Primera'
$ws.Cells.Item(49, 10).Value = 3000
$ws.Cells.Item(49, 11).Value = 1000
$ws.Cells.Item(49, 12).Value = 1000
$ws.Cells.Item(49, 13).Value = 1000
$ws.Cells.Item(49, 14).Value = '$/kilo'
$ws.Cells.Item(49, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(49, 16).Value = 1000
$ws.Cells.Item(49, 17).Value = 1
$ws.Cells.Item(49, 18).Value = 'Hortaliza'
